$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without letting Excel
# auto-convert numeric-looking strings (e.g. "318.00", "1.00") into numbers,
# and without leaving a stray style/number-format behind on the cell.
function Set-TextValue([string]$cellRef, [string]$value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Refresh prices / 1h volume percentages, and swap the FirstDigitalUSD / Celestia rows
# (rows 33 and 34), matching the upstream data refresh commit.
Set-TextValue "D2" "48.218.61"
Set-TextValue "E2" "  +0.04%  "
Set-TextValue "D3" "2.499.92"
Set-TextValue "E3" "  -1.24%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "318.00"
Set-TextValue "E5" "  -1.95%  "
Set-TextValue "D6" "106.11"
Set-TextValue "E6" "  -2.73%  "
Set-TextValue "E7" "  -1.67%  "
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  -0.03%  "
Set-TextValue "E9" "  -3.36%  "
Set-TextValue "D10" "39.03"
Set-TextValue "E10" "  -4.27%  "
Set-TextValue "D11" "20.28"
Set-TextValue "E11" "  -1.19%  "
Set-TextValue "D12" "0.0803"
Set-TextValue "E12" "  -3.03%  "
Set-TextValue "E13" "  +0.29%  "
Set-TextValue "D14" "7.11"
Set-TextValue "E14" "  -2.69%  "
Set-TextValue "D15" "2.893.11"
Set-TextValue "E15" "  -1.22%  "
Set-TextValue "D16" "2.514.19"
Set-TextValue "E16" "  -0.62%  "
Set-TextValue "D17" "0.830"
Set-TextValue "E17" "  -3.60%  "
Set-TextValue "D18" "48.114.42"
Set-TextValue "E18" "  +0.13%  "
Set-TextValue "D19" "2.99"
Set-TextValue "E19" "  +11.19%  "
Set-TextValue "D20" "12.84"
Set-TextValue "E20" "  -3.35%  "
Set-TextValue "E21" "  -1.07%  "
Set-TextValue "D22" "0.0₃0931"
Set-TextValue "E22" "  -2.21%  "
Set-TextValue "D23" "71.13"
Set-TextValue "D24" "267.66"
Set-TextValue "E24" "  -0.74%  "
Set-TextValue "E25" "  -2.59%  "
Set-TextValue "E26" "  +0.24%  "
Set-TextValue "D27" "25.80"
Set-TextValue "E27" "  -1.69%  "
Set-TextValue "E28" "  -0.63%  "
Set-TextValue "E29" "  -3.98%  "
Set-TextValue "E30" "  -3.00%  "
Set-TextValue "D31" "34.70"
Set-TextValue "E31" "  -2.92%  "
Set-TextValue "E32" "  -0.91%  "
Set-TextValue "B33" "FirstDigitalUSD"
Set-TextValue "C33" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  -0.09%  "
Set-TextValue "B34" "Celestia"
Set-TextValue "C34" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D34" "19.14"
Set-TextValue "E34" "  -3.81%  "
Set-TextValue "E35" "  -2.35%  "
Set-TextValue "E36" "  -2.71%  "
Set-TextValue "E37" "  -2.32%  "
Set-TextValue "D38" "4.60"
Set-TextValue "E38" "  -3.82%  "
Set-TextValue "D39" "2.89"
Set-TextValue "E39" "  -4.42%  "
Set-TextValue "D40" "123.18"
Set-TextValue "E40" "  +2.85%  "
Set-TextValue "E41" "  -1.91%  "
Set-TextValue "D42" "22.20"
Set-TextValue "E42" "  -0.83%  "
Set-TextValue "E43" "  +1.23%  "
Set-TextValue "D44" "0.0304"
Set-TextValue "E44" "  +0.74%  "
Set-TextValue "D45" "2.003.28"
Set-TextValue "E45" "  -0.58%  "
Set-TextValue "E46" "  +0.49%  "
Set-TextValue "D47" "1.92"
Set-TextValue "E47" "  +1.24%  "
Set-TextValue "E48" "  -2.58%  "
Set-TextValue "D49" "8.95"
Set-TextValue "E49" "  -2.47%  "
Set-TextValue "D50" "5.21"
Set-TextValue "E50" "  -1.15%  "
Set-TextValue "D51" "78.87"
Set-TextValue "E51" "  -1.25%  "
